# Auto-generated edit script: appends rows 15-21 to sheet "Artfynd"
# matching the unified diff (new fungi observation records).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 15 ---
$ws.Range("A15").Value = 111986331
$ws.Range("B15").Value = 90658
$ws.Range("C15").Value = 'Ovaliderad'
$ws.Range("D15").Value = 'NT'
$ws.Range("E15").Value = 4361
$ws.Range("F15").Value = 'Orange taggsvamp'
$ws.Range("G15").Value = 'Hydnellum aurantiacum'
$ws.Range("H15").Value = '(Batsch:Fr.) P.Karst.'
$ws.Range("I15").Value = "'"
$ws.Range("K15").Value = "'"
$ws.Range("P15").Value = 'Mellandammen (Mellandammen), Dlr'
$ws.Range("Q15").Value = 396461.0577280324
$ws.Range("R15").Value = 6849446.780411444
$ws.Range("S15").Value = 5
$ws.Range("T15").Value = 'Dalarna'
$ws.Range("U15").Value = 'Älvdalen'
$ws.Range("V15").Value = 'Dalarna'
$ws.Range("W15").Value = 'Särna'
$ws.Range("Y15").Value = "'2023-09-09"
$ws.Range("Z15").Value = '15:50'
$ws.Range("AA15").Value = "'2023-09-09"
$ws.Range("AB15").Value = '15:50'
$ws.Range("AD15").Value = $false
$ws.Range("AE15").Value = $false
$ws.Range("AG15").Value = $false
$ws.Range("AT15").Value = "'"
$ws.Range("AW15").Value = 'Bo karlstens'
$ws.Range("AX15").Value = 'Bo karlstens'
$ws.Range("AY15").Value = "'"

# --- Row 16 ---
$ws.Range("A16").Value = 111986181
$ws.Range("B16").Value = 90658
$ws.Range("C16").Value = 'Ovaliderad'
$ws.Range("D16").Value = 'NT'
$ws.Range("E16").Value = 4361
$ws.Range("F16").Value = 'Orange taggsvamp'
$ws.Range("G16").Value = 'Hydnellum aurantiacum'
$ws.Range("H16").Value = '(Batsch:Fr.) P.Karst.'
$ws.Range("I16").Value = "'"
$ws.Range("K16").Value = "'"
$ws.Range("P16").Value = 'Mellandammen (Mellandammen), Dlr'
$ws.Range("Q16").Value = 396464.1186182394
$ws.Range("R16").Value = 6849438.16119879
$ws.Range("S16").Value = 10
$ws.Range("T16").Value = 'Dalarna'
$ws.Range("U16").Value = 'Älvdalen'
$ws.Range("V16").Value = 'Dalarna'
$ws.Range("W16").Value = 'Särna'
$ws.Range("Y16").Value = "'2023-09-09"
$ws.Range("Z16").Value = '14:36'
$ws.Range("AA16").Value = "'2023-09-09"
$ws.Range("AB16").Value = '14:36'
$ws.Range("AD16").Value = $false
$ws.Range("AE16").Value = $false
$ws.Range("AG16").Value = $false
$ws.Range("AT16").Value = "'"
$ws.Range("AW16").Value = 'Bo karlstens'
$ws.Range("AX16").Value = 'Bo karlstens'
$ws.Range("AY16").Value = "'"

# --- Row 17 ---
$ws.Range("A17").Value = 111986477
$ws.Range("B17").Value = 90666
$ws.Range("C17").Value = 'Ovaliderad'
$ws.Range("D17").Value = 'LC'
$ws.Range("E17").Value = 4364
$ws.Range("F17").Value = 'Dropptaggsvamp'
$ws.Range("G17").Value = 'Hydnellum ferrugineum'
$ws.Range("H17").Value = '(Fr.:Fr.) P. Karst.'
$ws.Range("I17").Value = "'"
$ws.Range("K17").Value = "'"
$ws.Range("P17").Value = 'Mellandammen (Mellandammen), Dlr'
$ws.Range("Q17").Value = 396463.6812385211
$ws.Range("R17").Value = 6849392.223827818
$ws.Range("S17").Value = 5
$ws.Range("T17").Value = 'Dalarna'
$ws.Range("U17").Value = 'Älvdalen'
$ws.Range("V17").Value = 'Dalarna'
$ws.Range("W17").Value = 'Särna'
$ws.Range("Y17").Value = "'2023-09-09"
$ws.Range("Z17").Value = '15:58'
$ws.Range("AA17").Value = "'2023-09-09"
$ws.Range("AB17").Value = '15:58'
$ws.Range("AD17").Value = $false
$ws.Range("AE17").Value = $false
$ws.Range("AG17").Value = $false
$ws.Range("AT17").Value = "'"
$ws.Range("AW17").Value = 'Bo karlstens'
$ws.Range("AX17").Value = 'Bo karlstens'
$ws.Range("AY17").Value = "'"

# --- Row 18 ---
$ws.Range("A18").Value = 111986518
$ws.Range("B18").Value = 90678
$ws.Range("C18").Value = 'Ovaliderad'
$ws.Range("D18").Value = 'LC'
$ws.Range("E18").Value = 4366
$ws.Range("F18").Value = 'Skarp dropptaggsvamp'
$ws.Range("G18").Value = 'Hydnellum peckii'
$ws.Range("H18").Value = 'Banker'
$ws.Range("I18").Value = "'"
$ws.Range("K18").Value = "'"
$ws.Range("P18").Value = 'Mellandammen (Mellandammen), Dlr'
$ws.Range("Q18").Value = 396445.8145670656
$ws.Range("R18").Value = 6849381.867442117
$ws.Range("S18").Value = 5
$ws.Range("T18").Value = 'Dalarna'
$ws.Range("U18").Value = 'Älvdalen'
$ws.Range("V18").Value = 'Dalarna'
$ws.Range("W18").Value = 'Särna'
$ws.Range("Y18").Value = "'2023-09-09"
$ws.Range("Z18").Value = '16:06'
$ws.Range("AA18").Value = "'2023-09-09"
$ws.Range("AB18").Value = '16:06'
$ws.Range("AD18").Value = $false
$ws.Range("AE18").Value = $false
$ws.Range("AG18").Value = $false
$ws.Range("AT18").Value = "'"
$ws.Range("AW18").Value = 'Bo karlstens'
$ws.Range("AX18").Value = 'Bo karlstens'
$ws.Range("AY18").Value = "'"

# --- Row 19 ---
$ws.Range("A19").Value = 111986256
$ws.Range("B19").Value = 90660
$ws.Range("C19").Value = 'Ovaliderad'
$ws.Range("D19").Value = 'NT'
$ws.Range("E19").Value = 4362
$ws.Range("F19").Value = 'Blå taggsvamp'
$ws.Range("G19").Value = 'Hydnellum caeruleum'
$ws.Range("H19").Value = '(Hornem.) P.Karst.'
$ws.Range("I19").Value = "'"
$ws.Range("K19").Value = "'"
$ws.Range("P19").Value = 'Mellandammen (Mellandammen), Dlr'
$ws.Range("Q19").Value = 396466.3029379644
$ws.Range("R19").Value = 6849431.936995172
$ws.Range("S19").Value = 5
$ws.Range("T19").Value = 'Dalarna'
$ws.Range("U19").Value = 'Älvdalen'
$ws.Range("V19").Value = 'Dalarna'
$ws.Range("W19").Value = 'Särna'
$ws.Range("Y19").Value = "'2023-09-09"
$ws.Range("Z19").Value = '15:50'
$ws.Range("AA19").Value = "'2023-09-09"
$ws.Range("AB19").Value = '15:50'
$ws.Range("AD19").Value = $false
$ws.Range("AE19").Value = $false
$ws.Range("AG19").Value = $false
$ws.Range("AT19").Value = "'"
$ws.Range("AW19").Value = 'Bo karlstens'
$ws.Range("AX19").Value = 'Bo karlstens'
$ws.Range("AY19").Value = "'"

# --- Row 20 ---
$ws.Range("A20").Value = 111986412
$ws.Range("B20").Value = 90678
$ws.Range("C20").Value = 'Ovaliderad'
$ws.Range("D20").Value = 'LC'
$ws.Range("E20").Value = 4366
$ws.Range("F20").Value = 'Skarp dropptaggsvamp'
$ws.Range("G20").Value = 'Hydnellum peckii'
$ws.Range("H20").Value = 'Banker'
$ws.Range("I20").Value = "'"
$ws.Range("K20").Value = "'"
$ws.Range("P20").Value = 'Mellandammen (Mellandammen), Dlr'
$ws.Range("Q20").Value = 396473.4754867578
$ws.Range("R20").Value = 6849402.350115799
$ws.Range("S20").Value = 5
$ws.Range("T20").Value = 'Dalarna'
$ws.Range("U20").Value = 'Älvdalen'
$ws.Range("V20").Value = 'Dalarna'
$ws.Range("W20").Value = 'Särna'
$ws.Range("Y20").Value = "'2023-09-09"
$ws.Range("Z20").Value = '15:58'
$ws.Range("AA20").Value = "'2023-09-09"
$ws.Range("AB20").Value = '15:58'
$ws.Range("AD20").Value = $false
$ws.Range("AE20").Value = $false
$ws.Range("AG20").Value = $false
$ws.Range("AT20").Value = "'"
$ws.Range("AW20").Value = 'Bo karlstens'
$ws.Range("AX20").Value = 'Bo karlstens'
$ws.Range("AY20").Value = "'"

# --- Row 21 ---
$ws.Range("A21").Value = 111986397
$ws.Range("B21").Value = 90660
$ws.Range("C21").Value = 'Ovaliderad'
$ws.Range("D21").Value = 'NT'
$ws.Range("E21").Value = 4362
$ws.Range("F21").Value = 'Blå taggsvamp'
$ws.Range("G21").Value = 'Hydnellum caeruleum'
$ws.Range("H21").Value = '(Hornem.) P.Karst.'
$ws.Range("I21").Value = "'1"
$ws.Range("J21").Value = 'fruktkroppar'
$ws.Range("K21").Value = "'"
$ws.Range("P21").Value = 'Mellandammen (Mellandammen), Dlr'
$ws.Range("Q21").Value = 396473.4754867578
$ws.Range("R21").Value = 6849402.350115799
$ws.Range("S21").Value = 5
$ws.Range("T21").Value = 'Dalarna'
$ws.Range("U21").Value = 'Älvdalen'
$ws.Range("V21").Value = 'Dalarna'
$ws.Range("W21").Value = 'Särna'
$ws.Range("Y21").Value = "'2023-09-09"
$ws.Range("Z21").Value = '15:58'
$ws.Range("AA21").Value = "'2023-09-09"
$ws.Range("AB21").Value = '15:58'
$ws.Range("AD21").Value = $false
$ws.Range("AE21").Value = $false
$ws.Range("AG21").Value = $false
$ws.Range("AT21").Value = "'"
$ws.Range("AW21").Value = 'Bo karlstens'
$ws.Range("AX21").Value = 'Bo karlstens'
$ws.Range("AY21").Value = "'"

